$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A - shifts B:F left to A:E
$ws.Range("A1").EntireColumn.Delete()
